$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''48.546.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.83%  '

$ws.Range("D3").Value = '''2.611.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.67%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '''321.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.12%  '

$ws.Range("D6").Value = '''109.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("D7").Value = '''0.519'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.34%  '

$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").Value = '''0.537'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.85%  '

$ws.Range("D10").Value = '''39.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.13%  '

$ws.Range("E11").Value = '  -3.03%  '

$ws.Range("D12").Value = '''0.0806'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.96%  '

$ws.Range("E13").Value = '  +0.21%  '

$ws.Range("E14").Value = '  -0.98%  '

$ws.Range("D15").Value = '''3.032.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.11%  '

$ws.Range("D16").Value = '''2.595.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.27%  '

$ws.Range("D17").Value = '''0.857'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("D18").Value = '''48.539.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.63%  '

$ws.Range("D19").Value = '''2.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.57%  '

$ws.Range("D20").Value = '''12.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.63%  '

$ws.Range("D21").Value = '''6.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("D22").Value = '''0.0₃0936'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.55%  '

$ws.Range("D23").Value = '''268.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.57%  '

$ws.Range("D24").Value = '''68.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.69%  '

$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("D26").Value = '''25.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.66%  '

$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("D28").Value = '''9.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.74%  '

$ws.Range("D29").Value = '''2.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.71%  '

$ws.Range("D30").Value = '''0.136'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.56%  '

$ws.Range("D31").Value = '''34.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.12%  '

$ws.Range("D32").Value = '''49.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.83%  '

$ws.Range("D33").Value = '''5.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.60%  '

$ws.Range("E34").Value = '  -0.29%  '

$ws.Range("D35").Value = '''19.02'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.82%  '

$ws.Range("D36").Value = '''0.0790'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.03%  '

$ws.Range("D37").Value = '''4.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.31%  '

$ws.Range("E38").Value = '  +0.66%  '

$ws.Range("E39").Value = '  +4.84%  '

$ws.Range("D40").Value = '''125.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.03%  '

$ws.Range("D41").Value = '''22.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.35%  '

$ws.Range("E42").Value = '  -1.48%  '

$ws.Range("D43").Value = '''2.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.16%  '

$ws.Range("D44").Value = '''0.0312'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.83%  '

$ws.Range("D45").Value = '''2.051.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.03%  '

$ws.Range("E46").Value = '  -3.53%  '

$ws.Range("D47").Value = '''2.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.84%  '

$ws.Range("E48").Value = '  +2.07%  '

$ws.Range("D49").Value = '''8.86'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.31%  '

$ws.Range("D50").Value = '''58.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.32%  '

$ws.Range("D51").Value = '''5.13'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.02%  '
